# Remove the unused "CodeDark" custom layout (slideLayout2.xml) from the
# slide master. This is the only layout not referenced by any slide and
# matches the commit's removal of that layout from the deck's layout set.
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    if ($layout.Name -eq "CodeDark") {
        $layout.Delete()
        break
    }
}
